$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 11
$ws.Range("F4").Value = 947
$ws.Range("F5").Value = 1231
$ws.Range("F6").Value = 1676
$ws.Range("F7").Value = 896
$ws.Range("F9").Value = 2374
$ws.Range("F10").Value = 685
$ws.Range("F12").Value = 555
$ws.Range("F15").Value = 187
$ws.Range("F16").Value = 511
$ws.Range("F18").Value = 1223
$ws.Range("F21").Value = 2595
$ws.Range("F27").Value = 1734
$ws.Range("F29").Value = 524
$ws.Range("F34").Value = 75

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 7
$ws.Range("F11").Value = 61
$ws.Range("F20").Value = 58
$ws.Range("F24").Value = 192
$ws.Range("F33").Value = 22
$ws.Range("F35").Value = 62

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 497
$ws.Range("F7").Value = 155

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 497
$ws.Range("F8").Value = 947
$ws.Range("F9").Value = 1231
$ws.Range("F10").Value = 1676
$ws.Range("F14").Value = 896
$ws.Range("F16").Value = 2374
$ws.Range("F17").Value = 685
$ws.Range("F19").Value = 555
$ws.Range("F22").Value = 61
$ws.Range("F23").Value = 187
$ws.Range("F26").Value = 511
$ws.Range("F28").Value = 1223
$ws.Range("F32").Value = 2595
$ws.Range("F37").Value = 58
$ws.Range("F38").Value = 155
$ws.Range("F41").Value = 1734
$ws.Range("F42").Value = 524
$ws.Range("F48").Value = 75
$ws.Range("F49").Value = 62
